$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 and J1, copying the style from H1 (bold/bordered header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the numeric data for columns I and J, rows 2-6
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 6

$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 9
